# Added Test Data For Hungary/Russia/Finland Market
#
# Duplicates the "Denmark" template sheet three times to create new
# country sheets (Russia, Finland, Hungary) at the end of the workbook,
# fills in their market name / part-number cells, and leaves "Hungary"
# (the new last sheet) as the active / selected tab - matching the
# source diff.

$wb = $excel.ActiveWorkbook

# ---- Russia --------------------------------------------------------
$template = $wb.Worksheets.Item("Denmark")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Range("B4").Value = "NGC-2929/T2925"
$russia.Range("B2").Value = "Russia Market"

# ---- Finland --------------------------------------------------------
$template = $wb.Worksheets.Item("Denmark")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Range("B4").Value = "NGC-3130/T2887"
$finland.Range("B2").Value = "Finland Market"

# ---- Hungary ----------------------------------------------------------
$template = $wb.Worksheets.Item("Denmark")
$template.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Range("B4").Value = "NGC-3104/T2979"
$hungary.Range("B2").Value = "Hungary Market"

# Hungary becomes the active / selected tab (it was "Austria" before).
$hungary.Activate()
$hungary.Range("I17").Select()
